$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).Clear()
Write-Host "inserted and cleared"
